$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data in column A (script file names, appended below the
# existing table rows)
$ws.Range("A5").Value  = "SCRIPT/D16P31A/enter08.ssb"
$ws.Range("A6").Value  = "SCRIPT/D18P11A/enter06.ssb"
$ws.Range("A7").Value  = "SCRIPT/D19P11A/enter06.ssb"
$ws.Range("A8").Value  = "SCRIPT/D20P11A/enter06.ssb"
$ws.Range("A9").Value  = "SCRIPT/D22P11A/enter06.ssb"
$ws.Range("A10").Value = "SCRIPT/D23P11A/enter06.ssb"

# Column A got wider (to fit the longer file-name strings), column D
# got narrower
$ws.Columns(1).ColumnWidth = 29
$ws.Columns(4).ColumnWidth = 24.2

# Header row shrinks to a single line; the wrapped-text data rows go
# back to automatic (content-based) height now that the columns changed
$ws.Rows(1).RowHeight = 29.4
$ws.Rows(2).AutoFit()
$ws.Rows(3).AutoFit()
$ws.Rows(4).AutoFit()

# Active selection ends on B6
$ws.Range("B6").Select() | Out-Null
